# Apply cryptocurrency price/volume updates scraped on Mon Jun  5 19:16:53 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.566.28"
$ws.Range("E2").Value = "  -6.14%  "
$ws.Range("D3").Value = "1.807.65"
$ws.Range("E3").Value = "  -5.20%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "276.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -9.79%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5003"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -6.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3500"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -8.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.83"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06651"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -8.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8360"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07826"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.71%  "
$ws.Range("D14").Value = "1.812.16"
$ws.Range("E14").Value = "  +69.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.038"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.11"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.76%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007870"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.12%  "
$ws.Range("D21").Value = "25.645.00"
$ws.Range("E21").Value = "  -5.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.712"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.64%  "
$ws.Range("D23").Value = "2.035.08"
$ws.Range("E23").Value = "  +69.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.976"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.047"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.27%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.107"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -8.51%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.662"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.00%  "
$ws.Range("E29").Value = "  -8.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "108.56"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -7.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.288"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -11.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.199"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -10.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08808"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04791"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7328"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -11.64%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.120"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.15%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.850"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.000"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.035"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01859"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5188"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -12.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.299"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -14.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9588"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -11.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.50"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.164"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.037"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -14.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4582"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -9.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1376"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -10.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.180"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.53%  "
